$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-capitalize the header row: "id"/"name" -> "Id"/"Name"
$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "Name"

# Append a new data row: Id=3, Name="Eve"
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Eve"

# Match the formatting (wrap text / vertical-center) used by the other data rows
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to match the saved view state
$ws.Range("C5").Select()
